$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 582
$ws.Range("F3").Value = 262
$ws.Range("F5").Value = 743
$ws.Range("F6").Value = 377
$ws.Range("F8").Value = 157
$ws.Range("F9").Value = 242
$ws.Range("F10").Value = 226
$ws.Range("F11").Value = 6031
$ws.Range("F12").Value = 58
$ws.Range("F13").Value = 49
$ws.Range("F14").Value = 497
$ws.Range("F16").Value = 548
$ws.Range("F17").Value = 362
$ws.Range("F21").Value = 711
$ws.Range("F22").Value = 149
$ws.Range("F24").Value = 316
$ws.Range("F25").Value = 1020
$ws.Range("F27").Value = 1825
$ws.Range("F28").Value = 490

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 268
$ws.Range("F5").Value = 271
$ws.Range("F6").Value = 300

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 248

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 248
$ws.Range("F3").Value = 582
$ws.Range("F4").Value = 262
$ws.Range("F6").Value = 743
$ws.Range("F8").Value = 377
$ws.Range("F10").Value = 157
$ws.Range("F11").Value = 242
$ws.Range("F12").Value = 226
$ws.Range("F13").Value = 6031
$ws.Range("F14").Value = 58
$ws.Range("F15").Value = 49
$ws.Range("F16").Value = 268
$ws.Range("F17").Value = 497
$ws.Range("F19").Value = 548
$ws.Range("F20").Value = 362
$ws.Range("F25").Value = 271
$ws.Range("F26").Value = 300
$ws.Range("F28").Value = 711
$ws.Range("F32").Value = 149
$ws.Range("F34").Value = 316
$ws.Range("F35").Value = 1020
$ws.Range("F37").Value = 1825
$ws.Range("F38").Value = 490
